$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B11").Value = 200
$ws.Range("C11").Value = 5000
$ws.Range("B65").Value = 1

$ws.Range("B65").Select()
